$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2023 year column (K) with new data
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1135.7
$ws.Range("K5").Value = 791.1
$ws.Range("K6").Value = 1360

# Copy style from column J to column K for rows 3-6
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
